# Sync automático del tracker - actualiza el resultado de predicciones
# completadas que estaban en estado "Pending".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Predictions")

$timestamp = "2025-08-23 02:25:43"

# Filas a sincronizar: numero de fila -> (Result, Profit)
$updates = @(
    @{ Row = 31; Result = "Home Win"; Profit = -4.3 },
    @{ Row = 33; Result = "Draw";     Profit = -4.1 },
    @{ Row = 40; Result = "Draw";     Profit = -2.7 },
    @{ Row = 41; Result = "Home Win"; Profit = -6.2 }
)

foreach ($u in $updates) {
    $r = $u.Row

    $ws.Range("L$r").Value = "Completed"       # Status
    $ws.Range("M$r").Value = $u.Result          # Result
    $ws.Range("N$r").Value = "Fallo"            # Resultado_Real
    $ws.Range("O$r").Value = $u.Profit          # Profit
    $ws.Range("P$r").Value = -100               # ROI
    $ws.Range("Q$r").Value = $timestamp         # Enviado
}
